# Weekly update: insert a new price record for "Piña / Caramelo / Tercera"
# (Macroferia Regional de Talca) as row 232, shifting the existing rows
# 232-295 down to 233-296.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 232 (pushes rows 232.. down by one).
$ws.Rows.Item(232).Insert()

# Populate the new row 232 with the new weekly data point.
$ws.Cells.Item(232, 1).Value = 5
$ws.Cells.Item(232, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(232, 3).Value = "Maule"
$ws.Cells.Item(232, 4).Value = 44841
$ws.Cells.Item(232, 5).Value = 7
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100108
$ws.Cells.Item(232, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(232, 9).Value = 100108005
$ws.Cells.Item(232, 10).Value = "Piña"
$ws.Cells.Item(232, 11).Value = "Caramelo"
$ws.Cells.Item(232, 12).Value = "Tercera"
$ws.Cells.Item(232, 13).Value = 220
$ws.Cells.Item(232, 14).Value = 21000
$ws.Cells.Item(232, 15).Value = 21000
$ws.Cells.Item(232, 16).Value = 21000
$ws.Cells.Item(232, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(232, 18).Value = "Ecuador"
$ws.Cells.Item(232, 19).Value = 1312
$ws.Cells.Item(232, 20).Value = 16
